$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Student ID values in rows 2-7 (use a leading apostrophe so the
# numeric-looking IDs stay stored as text, matching the source data)
$ws.Range("A2").Value = "'201795"
$ws.Range("A3").Value = "'211177"
$ws.Range("A4").Value = "'190846"
$ws.Range("A5").Value = "'201563"
$ws.Range("A6").Value = "'211131"
$ws.Range("A7").Value = "'200359"

# Remove the last two rows (rows 8 and 9) entirely
$ws.Rows("8:9").Delete()
